$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Translate the Swedish "Rad X, Kol Y" labels to English "Row X, Col Y" ---
# The sheet is a 7-row x 6-col grid of labeled cells (A1:F7).
for ($r = 1; $r -le 7; $r++) {
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = "Row " + $r.ToString() + ", Col " + $c.ToString()
    }
}

# --- Give rows 1-9 an explicit 15.75pt height (this also creates rows 8 & 9) ---
for ($r = 1; $r -le 9; $r++) {
    $ws.Rows.Item($r).RowHeight = 15.75
}

# --- Give columns C:F an explicit width matching the sheet's default column width ---
$ws.Columns("C:F").ColumnWidth = 11.8
